# Scheduled market-data refresh: update currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the per-job sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) with freshly polled Universalis prices.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1048.7333
$ws.Range("I19").Value = 981.2222
$ws.Range("J19").Value = 1150
$ws.Range("K19").Value = 981.2222
$ws.Range("L19").Value = 1150
$ws.Range("M19").Value = -806.2222
$ws.Range("N19").Value = -1500
$ws.Range("H46").Value = 27530
$ws.Range("J46").Value = 31703.846
$ws.Range("L46").Value = 95111.538
$ws.Range("N46").Value = -95349.538
$ws.Range("H60").Value = 27530
$ws.Range("J60").Value = 31703.846
$ws.Range("L60").Value = 95111.538
$ws.Range("N60").Value = -96079.538
$ws.Range("H75").Value = 55314
$ws.Range("J75").Value = 55314
$ws.Range("L75").Value = 55314
$ws.Range("N75").Value = -57186
$ws.Range("H78").Value = 55314
$ws.Range("J78").Value = 55314
$ws.Range("L78").Value = 165942
$ws.Range("N78").Value = -175302
$ws.Range("H100").Value = 37040696
$ws.Range("I100").Value = 66668656
$ws.Range("J100").Value = 5750
$ws.Range("K100").Value = 66668656
$ws.Range("L100").Value = 5750
$ws.Range("M100").Value = -66668115
$ws.Range("N100").Value = -6832
$ws.Range("H112").Value = 1152.5
$ws.Range("J112").Value = 1345.8334
$ws.Range("L112").Value = 4037.5002
$ws.Range("N112").Value = -6253.5002
$ws.Range("H127").Value = 1376.0344
$ws.Range("I127").Value = 637.3333
$ws.Range("J127").Value = 1568.7391
$ws.Range("K127").Value = 1911.9999
$ws.Range("L127").Value = 4706.2173
$ws.Range("M127").Value = 3048.0001
$ws.Range("N127").Value = -14626.2173
$ws.Range("H129").Value = 1072.4222
$ws.Range("I129").Value = 496.8
$ws.Range("J129").Value = 1144.375
$ws.Range("K129").Value = 1490.4
$ws.Range("L129").Value = 3433.125
$ws.Range("M129").Value = 3509.6
$ws.Range("N129").Value = -13433.125
$ws.Range("H138").Value = 4828.107
$ws.Range("J138").Value = 4828.107
$ws.Range("L138").Value = 14484.321
$ws.Range("N138").Value = -24764.321

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1162.5555
$ws.Range("I61").Value = 993.8333
$ws.Range("K61").Value = 993.8333
$ws.Range("M61").Value = -781.8333
$ws.Range("H68").Value = 33549.5
$ws.Range("J68").Value = 33549.5
$ws.Range("L68").Value = 33549.5
$ws.Range("N68").Value = -35171.5
$ws.Range("H71").Value = 33549.5
$ws.Range("J71").Value = 33549.5
$ws.Range("L71").Value = 100648.5
$ws.Range("N71").Value = -108760.5
$ws.Range("H136").Value = 1162.5555
$ws.Range("I136").Value = 993.8333
$ws.Range("K136").Value = 2981.4999
$ws.Range("M136").Value = -431.4998999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 500
$ws.Range("I76").Value = 500
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 500
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("M76").Value = -185
$ws.Range("H79").Value = 500
$ws.Range("I79").Value = 500
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 500
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("M79").Value = 592
$ws.Range("H99").Value = 2280
$ws.Range("I99").Value = 1906.6666
$ws.Range("J99").Value = 2381.818
$ws.Range("K99").Value = 1906.6666
$ws.Range("L99").Value = 2381.818
$ws.Range("M99").Value = -408.6666
$ws.Range("N99").Value = -5377.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 16050000
$ws.Range("I99").Value = 32000000
$ws.Range("J99").Value = 100000
$ws.Range("K99").Value = 32000000
$ws.Range("L99").Value = 100000
$ws.Range("M99").Value = -31998502
$ws.Range("N99").Value = -102996
$ws.Range("H126").Value = 16050000
$ws.Range("I126").Value = 32000000
$ws.Range("J126").Value = 100000
$ws.Range("K126").Value = 96000000
$ws.Range("L126").Value = 300000
$ws.Range("M126").Value = -95997530
$ws.Range("N126").Value = -304940

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1200.7
$ws.Range("I113").Value = 3001.5
$ws.Range("J113").Value = 750.5
$ws.Range("K113").Value = 9004.5
$ws.Range("L113").Value = 2251.5
$ws.Range("M113").Value = -6834.5
$ws.Range("N113").Value = -6591.5
$ws.Range("H131").Value = 4034.6365
$ws.Range("J131").Value = 5031.269
$ws.Range("L131").Value = 15093.807
$ws.Range("N131").Value = -25173.807

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2534.25
$ws.Range("I102").Value = 2400.3333
$ws.Range("J102").Value = 2706.4285
$ws.Range("K102").Value = 2400.3333
$ws.Range("L102").Value = 2706.4285
$ws.Range("M102").Value = -778.3332999999998
$ws.Range("N102").Value = -5950.4285
$ws.Range("H136").Value = 16339.571
$ws.Range("J136").Value = 16339.571
$ws.Range("L136").Value = 49018.713
$ws.Range("N136").Value = -54118.713

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4335.6665
$ws.Range("I122").Value = 4001.3333
$ws.Range("J122").Value = 4670
$ws.Range("K122").Value = 12003.9999
$ws.Range("L122").Value = 14010
$ws.Range("M122").Value = -9553.999899999999
$ws.Range("N122").Value = -18910
$ws.Range("H136").Value = 4872.364
$ws.Range("I136").Value = 5209.6
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 15628.8
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -13078.8
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 20000
$ws.Range("J80").Value = 20000
$ws.Range("L80").Value = 20000
$ws.Range("N80").Value = -21996
$ws.Range("H83").Value = 20000
$ws.Range("J83").Value = 20000
$ws.Range("L83").Value = 60000
$ws.Range("N83").Value = -69984
$ws.Range("H122").Value = 87398.164
$ws.Range("I122").Value = 114866.89
$ws.Range("J122").Value = 4992
$ws.Range("K122").Value = 344600.67
$ws.Range("L122").Value = 14976
$ws.Range("M122").Value = -342150.67
$ws.Range("N122").Value = -19876
$ws.Range("H126").Value = 44509.305
$ws.Range("I126").Value = 48591.145
$ws.Range("J126").Value = 1650
$ws.Range("K126").Value = 145773.435
$ws.Range("L126").Value = 4950
$ws.Range("M126").Value = -143303.435
$ws.Range("N126").Value = -9890
